# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# The "K" column (column G) previously held strike-count totals; it is
# regenerated here to hold actual strikeout (K) counts per game.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 5
    3  = 4
    4  = 5
    5  = 9
    6  = 5
    7  = 6
    8  = 8
    9  = 6
    10 = 10
    11 = 4
    12 = 4
    13 = 7
    14 = 6
    15 = 3
    16 = 3
    17 = 8
    18 = 3
    19 = 1
    20 = 5
    21 = 6
    22 = 7
    23 = 4
    24 = 6
    25 = 5
    26 = 4
    27 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
